$wb = $excel.ActiveWorkbook

# Original sheet (Sheet1 / sheet1.xml) becomes "Sheet0" - keep its data,
# but it is no longer the active/selected tab; selection moves to A4.
$sheet0 = $wb.Worksheets.Item(1)
$sheet0.Name = "Sheet0"

# New sheet inserted after Sheet0 -> becomes "Sheet1" (sheet2.xml). This is
# a re-creation of the original data with row 2 missing, and becomes the
# active tab.
$sheet1 = $wb.Worksheets.Add($null, $sheet0)
$sheet1.Name = "Sheet1"
$sheet1.Range("A1").Value = 2
$sheet1.Range("A3").Value = 3

# New empty sheet inserted after Sheet1 -> becomes "Sheet6" (sheet3.xml).
$sheet6 = $wb.Worksheets.Add($null, $sheet1)
$sheet6.Name = "Sheet6"

# Sheet0's (sheet1.xml) selection moves to A4. Select it while Sheet0 is
# still active so the selection is recorded against that sheet.
$sheet0.Activate()
$sheet0.Range("A4").Select()

# Sheet1 (sheet2.xml) is the active tab (index 1 / activeTab="1"), with its
# sheetView tabSelected and default A1 selection. Activate it last so it
# ends up as the saved active sheet.
$sheet1.Activate()
